# Update date heading
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-23 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-24 Tuesday", 2)

# Update the division problems in the table. Cell text is addressed by
# row/column rather than by Find/Replace because some problems (e.g. the
# original "22÷5=") occur more than once with different replacements.
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "40÷4="
$table.Cell(1, 2).Range.Text = "36÷2="
$table.Cell(1, 3).Range.Text = "52÷7="
$table.Cell(1, 4).Range.Text = "96÷9="
$table.Cell(1, 5).Range.Text = "37÷4="

$table.Cell(5, 1).Range.Text = "21÷2="
$table.Cell(5, 2).Range.Text = "61÷3="
$table.Cell(5, 3).Range.Text = "80÷4="
$table.Cell(5, 4).Range.Text = "94÷4="
$table.Cell(5, 5).Range.Text = "69÷8="

$table.Cell(9, 1).Range.Text = "97÷7="
$table.Cell(9, 2).Range.Text = "18÷9="
$table.Cell(9, 3).Range.Text = "54÷2="
$table.Cell(9, 4).Range.Text = "25÷8="
$table.Cell(9, 5).Range.Text = "71÷2="

$table.Cell(13, 1).Range.Text = "91÷3="
$table.Cell(13, 2).Range.Text = "53÷5="
$table.Cell(13, 3).Range.Text = "65÷5="
$table.Cell(13, 4).Range.Text = "37÷9="
$table.Cell(13, 5).Range.Text = "19÷3="

$table.Cell(17, 1).Range.Text = "51÷7="
$table.Cell(17, 2).Range.Text = "27÷5="
$table.Cell(17, 3).Range.Text = "49÷2="
$table.Cell(17, 4).Range.Text = "10÷4="
$table.Cell(17, 5).Range.Text = "90÷3="
